$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.724.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "'1.574.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'213.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'44.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("D9").Value = "'24.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "'1.799.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'1.581.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "'28.723.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "'3.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "'62.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "'230.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -4.85%  "
$ws.Range("D24").Value = "'9.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  +8.58%  "
$ws.Range("D26").Value = "'152.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "'1.394.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  +3.31%  "
$ws.Range("E37").Value = "  -2.76%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("E39").Value = "  +3.16%  "
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("D41").Value = "'0.527"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'0.794"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").Value = "'0.0472"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "'0.966"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").Value = "'63.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").Value = "'1.711.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "'86.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("E51").Value = "  -0.93%  "
